$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")
Write-Host $ws.Name
Write-Host $ws.Range("A28").Value
